$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The betting-odds rows below were reshuffled: each rows B:AC payload (id,
# Div, Div Original Name, Date, HomeTeam, AwayTeam, odds, etc.) moved to a
# different row while the row index in column A stayed put. We snapshot the
# current B:AC values of every row in a cycle first, then write them back in
# the new order so a row that is both a source and a destination is not
# clobbered before it has been read.

# Cycle: 368, 369, 371, 370
$row368 = $ws.Range("B368:AC368").Value()
$row369 = $ws.Range("B369:AC369").Value()
$row371 = $ws.Range("B371:AC371").Value()
$row370 = $ws.Range("B370:AC370").Value()
$ws.Range("B368:AC368").Value = $row369
$ws.Range("B369:AC369").Value = $row371
$ws.Range("B371:AC371").Value = $row370
$ws.Range("B370:AC370").Value = $row368

# Cycle: 373, 375, 374
$row373 = $ws.Range("B373:AC373").Value()
$row375 = $ws.Range("B375:AC375").Value()
$row374 = $ws.Range("B374:AC374").Value()
$ws.Range("B373:AC373").Value = $row375
$ws.Range("B375:AC375").Value = $row374
$ws.Range("B374:AC374").Value = $row373

# Cycle: 377, 378, 379, 380, 381
$row377 = $ws.Range("B377:AC377").Value()
$row378 = $ws.Range("B378:AC378").Value()
$row379 = $ws.Range("B379:AC379").Value()
$row380 = $ws.Range("B380:AC380").Value()
$row381 = $ws.Range("B381:AC381").Value()
$ws.Range("B377:AC377").Value = $row378
$ws.Range("B378:AC378").Value = $row379
$ws.Range("B379:AC379").Value = $row380
$ws.Range("B380:AC380").Value = $row381
$ws.Range("B381:AC381").Value = $row377
